# Add 2022-Q3 data:
#  - Insert a new worksheet "2022-Q3" right before "2022-Q2" with the
#    quarterly fund-holdings detail.
#  - Insert a new row at the top of the summary ("总计") sheet with the
#    2022-Q3 totals, shifting the existing history rows down by one.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) New detail sheet "2022-Q3", positioned before "2022-Q2".
#    Duplicate the "2022-Q2" sheet so fonts/borders/column layout are
#    preserved exactly, then overwrite the data with the Q3 figures.
# ---------------------------------------------------------------------
$templateSheet = $wb.Worksheets.Item("2022-Q2")
$templateSheet.Copy($templateSheet)
$q3 = $wb.Worksheets.Item("2022-Q2 (2)")
$q3.Name = "2022-Q3"

# The template has 3 data rows (rows 2-4); Q3 only needs 1 (row 2), so
# drop the extra two rows.
$q3.Rows.Item(3).Delete()
$q3.Rows.Item(3).Delete()

# Header row (row 1) is identical across quarters, so it's already
# correct from the copy. Only the data row changes.
$q3.Range("A2").Value = 0
$q3.Range("B2").Value = "'014887"
$q3.Range("C2").Value = "招商安福1年定期开放债券"
$q3.Range("D2").Value = "'17.28"
$q3.Range("E2").Value = "'29.21"
$q3.Range("F2").Value = "'1.08"
$q3.Range("G2").Value = "'0.1866"
$q3.Range("H2").Value = 10

# ---------------------------------------------------------------------
# 2) Update the "总计" (summary) sheet: push the existing quarters down
#    one row and write the new 2022-Q3 row at the top (row 2).
# ---------------------------------------------------------------------
$total = $wb.Worksheets.Item("总计")

for ($r = 7; $r -ge 2; $r--) {
    $dest = $r + 1
    $total.Cells.Item($dest, 1).Value = $total.Cells.Item($r, 1).Value()
    $total.Cells.Item($dest, 2).Value = $total.Cells.Item($r, 2).Value()
    $total.Cells.Item($dest, 3).Value = $total.Cells.Item($r, 3).Value()
    $total.Cells.Item($dest, 4).Value = $total.Cells.Item($r, 4).Value()
}

# Row 8 is brand-new: give its A-column cell the same formatting as the
# other A-column entries (bold/centered/bordered header style).
$total.Cells.Item(2, 1).Copy()
$total.Cells.Item(8, 1).PasteSpecial(-4122)

$total.Range("A2").Value = 0
$total.Range("B2").Value = "2022-Q3"
$total.Range("C2").Value = 1
$total.Range("D2").Value = 0.19
